# Update 'F' column (想去人数 / interest count) values per the scraped data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 2378
$ws.Cells.Item(3, 6).Value = 580
$ws.Cells.Item(4, 6).Value = 211
$ws.Cells.Item(5, 6).Value = 368
$ws.Cells.Item(6, 6).Value = 368
$ws.Cells.Item(7, 6).Value = 613
$ws.Cells.Item(9, 6).Value = 820
$ws.Cells.Item(10, 6).Value = 544
$ws.Cells.Item(11, 6).Value = 845
$ws.Cells.Item(12, 6).Value = 392
$ws.Cells.Item(13, 6).Value = 103
$ws.Cells.Item(14, 6).Value = 406
$ws.Cells.Item(15, 6).Value = 24
$ws.Cells.Item(16, 6).Value = 1040
$ws.Cells.Item(17, 6).Value = 21963
$ws.Cells.Item(18, 6).Value = 1001
$ws.Cells.Item(19, 6).Value = 100
$ws.Cells.Item(20, 6).Value = 291
$ws.Cells.Item(21, 6).Value = 320
$ws.Cells.Item(23, 6).Value = 182
$ws.Cells.Item(25, 6).Value = 23
$ws.Cells.Item(26, 6).Value = 268
$ws.Cells.Item(28, 6).Value = 374
$ws.Cells.Item(29, 6).Value = 166

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 37
$ws.Cells.Item(5, 6).Value = 91
$ws.Cells.Item(6, 6).Value = 212
$ws.Cells.Item(7, 6).Value = 233
$ws.Cells.Item(8, 6).Value = 3472
$ws.Cells.Item(10, 6).Value = 116
$ws.Cells.Item(16, 6).Value = 3992

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 273
$ws.Cells.Item(3, 6).Value = 126
$ws.Cells.Item(4, 6).Value = 649

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 273
$ws.Cells.Item(3, 6).Value = 126
$ws.Cells.Item(4, 6).Value = 37
$ws.Cells.Item(5, 6).Value = 2378
$ws.Cells.Item(6, 6).Value = 649
$ws.Cells.Item(7, 6).Value = 580
$ws.Cells.Item(8, 6).Value = 211
$ws.Cells.Item(9, 6).Value = 368
$ws.Cells.Item(10, 6).Value = 368
$ws.Cells.Item(11, 6).Value = 613
$ws.Cells.Item(15, 6).Value = 91
$ws.Cells.Item(16, 6).Value = 212
$ws.Cells.Item(18, 6).Value = 820
$ws.Cells.Item(19, 6).Value = 544
$ws.Cells.Item(20, 6).Value = 845
$ws.Cells.Item(21, 6).Value = 392
$ws.Cells.Item(22, 6).Value = 103
$ws.Cells.Item(23, 6).Value = 406
$ws.Cells.Item(24, 6).Value = 24
$ws.Cells.Item(25, 6).Value = 1040
$ws.Cells.Item(26, 6).Value = 21964
$ws.Cells.Item(27, 6).Value = 233
$ws.Cells.Item(28, 6).Value = 3472
$ws.Cells.Item(30, 6).Value = 116
$ws.Cells.Item(32, 6).Value = 1001
$ws.Cells.Item(33, 6).Value = 100
$ws.Cells.Item(34, 6).Value = 291
$ws.Cells.Item(37, 6).Value = 320
$ws.Cells.Item(39, 6).Value = 182
$ws.Cells.Item(41, 6).Value = 23
$ws.Cells.Item(44, 6).Value = 268
$ws.Cells.Item(46, 6).Value = 374
$ws.Cells.Item(47, 6).Value = 166
$ws.Cells.Item(48, 6).Value = 3992
